$wb = $excel.ActiveWorkbook

# --- New handoff file identifiers (this run) ---
$newGuid   = "84863cc1-188a-4b60-b18d-6e148d2b4d34"
$newMd     = "$newGuid.md"
$newHash   = "61aa76ce3402e706981373dde77be4e7046e540a"
$zhXlf     = "$newGuid.$newHash.zh-cn.xlf"
$deXlf     = "$newGuid.$newHash.de-de.xlf"

$mdUrl    = "https://github.com/OpenLocalizationTest/oltest/blob/0abfdd0cba2cae886fb8e4388caa1c5d623b3e4d/e2e/$newMd"
$zhXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/04d5e98b4ea334717a3e4c0b03eda8f5eed97f53/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$zhXlf"
$deXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/871ad7d9a57c4100f021a6afffd47f2d6c54645c/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$deXlf"

# =====================================================================
# Sheet "Overview" -- add row 3 for the new handoff
# =====================================================================
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = "2016-27-18 08:27:22"
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), $mdUrl, "", "", $newMd)

# =====================================================================
# Sheet "zh-cn" -- add row 3 for the new handoff
# =====================================================================
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C3").Value = "Ready for handoff"
$wsZh.Range("E3").Value = "2016-03-18 08:27:19"
$wsZh.Range("H3").Value = "0001-01-01 00:00:00"
$wsZh.Range("I3").Value = "Include"

$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $mdUrl, "", "", $newMd)
$wsZh.Hyperlinks.Add($wsZh.Range("B3"), $mdUrl, "", "", ".md")
$wsZh.Hyperlinks.Add($wsZh.Range("D3"), $zhXlfUrl, "", "", $zhXlf)

# =====================================================================
# Sheet "de-de" -- add row 3 for the new handoff
# =====================================================================
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C3").Value = "Ready for handoff"
$wsDe.Range("E3").Value = "2016-03-18 08:27:22"
$wsDe.Range("H3").Value = "0001-01-01 00:00:00"
$wsDe.Range("I3").Value = "Include"

$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $mdUrl, "", "", $newMd)
$wsDe.Hyperlinks.Add($wsDe.Range("B3"), $mdUrl, "", "", ".md")
$wsDe.Hyperlinks.Add($wsDe.Range("D3"), $deXlfUrl, "", "", $deXlf)
